{"js": "// Replace the old two-digit \u00f7 one-digit division problems with the new\n// ones, per the commit's regenerated output. Each \"old\" string is unique\n// in the document, so a direct search+replace per pair is unambiguous.\nconst replacements = [\n  [\"19\u00f74=4, 3\", \"62\u00f79=6, 8\"],\n  [\"52\u00f79=5, 7\", \"20\u00f79=2, 2\"],\n  [\"63\u00f78=7, 7\", \"28\u00f73=9, 1\"],\n  [\"97\u00f73=32, 1\", \"46\u00f76=7, 4\"],\n  [\"94\u00f77=13, 3\", \"26\u00f73=8, 2\"],\n  [\"10\u00f77=1, 3\", \"74\u00f72=37, 0\"],\n  [\"32\u00f78=4, 0\", \"36\u00f72=18, 0\"],\n  [\"45\u00f74=11, 1\", \"29\u00f77=4, 1\"],\n  [\"50\u00f73=16, 2\", \"82\u00f72=41, 0\"],\n  [\"73\u00f78=9, 1\", \"59\u00f78=7, 3\"],\n  [\"39\u00f75=7, 4\", \"46\u00f75=9, 1\"],\n  [\"11\u00f76=1, 5\", \"26\u00f75=5, 1\"],\n  [\"49\u00f77=7, 0\", \"33\u00f78=4, 1\"],\n  [\"31\u00f79=3, 4\", \"69\u00f72=34, 1\"],\n  [\"23\u00f78=2, 7\", \"78\u00f75=15, 3\"],\n  [\"66\u00f73=22, 0\", \"53\u00f74=13, 1\"],\n  [\"38\u00f77=5, 3\", \"28\u00f78=3, 4\"],\n  [\"81\u00f78=10, 1\", \"41\u00f77=5, 6\"],\n  [\"52\u00f73=17, 1\", \"20\u00f78=2, 4\"],\n  [\"67\u00f79=7, 4\", \"35\u00f75=7, 0\"],\n  [\"89\u00f75=17, 4\", \"74\u00f72=37, 0\"],\n  [\"51\u00f79=5, 6\", \"35\u00f73=11, 2\"],\n  [\"34\u00f73=11, 1\", \"33\u00f72=16, 1\"],\n  [\"68\u00f75=13, 3\", \"47\u00f76=7, 5\"],\n  [\"78\u00f72=39, 0\", \"49\u00f78=6, 1\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  if (found.items.length === 0) {\n    throw new Error(`Could not find text to replace: \"${oldText}\"`);\n  }\n\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the old two-digit / one-digit division problems with the new\n# ones, per the commit's regenerated output. Each \"old\" string is unique\n# in the document, so Find/Replace-All per pair is unambiguous (one hit\n# each).\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"19\u00f74=4, 3\", \"62\u00f79=6, 8\"),\n  @(\"52\u00f79=5, 7\", \"20\u00f79=2, 2\"),\n  @(\"63\u00f78=7, 7\", \"28\u00f73=9, 1\"),\n  @(\"97\u00f73=32, 1\", \"46\u00f76=7, 4\"),\n  @(\"94\u00f77=13, 3\", \"26\u00f73=8, 2\"),\n  @(\"10\u00f77=1, 3\", \"74\u00f72=37, 0\"),\n  @(\"32\u00f78=4, 0\", \"36\u00f72=18, 0\"),\n  @(\"45\u00f74=11, 1\", \"29\u00f77=4, 1\"),\n  @(\"50\u00f73=16, 2\", \"82\u00f72=41, 0\"),\n  @(\"73\u00f78=9, 1\", \"59\u00f78=7, 3\"),\n  @(\"39\u00f75=7, 4\", \"46\u00f75=9, 1\"),\n  @(\"11\u00f76=1, 5\", \"26\u00f75=5, 1\"),\n  @(\"49\u00f77=7, 0\", \"33\u00f78=4, 1\"),\n  @(\"31\u00f79=3, 4\", \"69\u00f72=34, 1\"),\n  @(\"23\u00f78=2, 7\", \"78\u00f75=15, 3\"),\n  @(\"66\u00f73=22, 0\", \"53\u00f74=13, 1\"),\n  @(\"38\u00f77=5, 3\", \"28\u00f78=3, 4\"),\n  @(\"81\u00f78=10, 1\", \"41\u00f77=5, 6\"),\n  @(\"52\u00f73=17, 1\", \"20\u00f78=2, 4\"),\n  @(\"67\u00f79=7, 4\", \"35\u00f75=7, 0\"),\n  @(\"89\u00f75=17, 4\", \"74\u00f72=37, 0\"),\n  @(\"51\u00f79=5, 6\", \"35\u00f73=11, 2\"),\n  @(\"34\u00f73=11, 1\", \"33\u00f72=16, 1\"),\n  @(\"68\u00f75=13, 3\", \"47\u00f76=7, 5\"),\n  @(\"78\u00f72=39, 0\", \"49\u00f78=6, 1\")\n)\n\nforeach ($pair in $pairs) {\n  $old = $pair[0]\n  $new = $pair[1]\n\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $old\n  $find.Replacement.Text = $new\n  $matchCase = $true\n  $matchWholeWord = $false\n  $matchWildcards = $false\n  $matchSoundsLike = $false\n  $matchAllWordForms = $false\n  $forward = $true\n  $wrap = 1\n  $format = $false\n  $replace = 2\n\n  $found = $find.Execute($old, $matchCase, $matchWholeWord, $matchWildcards, $matchSoundsLike, $matchAllWordForms, $forward, $wrap, $format, $new, $replace)\n\n  if (-not $found) {\n    throw \"Could not find text to replace: '$old'\"\n  }\n}\n"}
